$d = $word.ActiveDocument

# NOTE: all search/replace strings use single quotes so PowerShell does not try to
# expand ${...} template placeholders as variables.

# --- Sentence 1 (paragraph starting with ${tenCoQuanDuocGiaoNhiemVuXM} bao cao...) ---
# Originally split across 3 runs:
#   "...nêu trên vớ" | "i ${tenLanhDao}" | " trước ngày... tháng ... năm ..."
# Target: merge into a single run with the full text unchanged.
$old1 = 'nêu trên vớ' + 'i ${tenLanhDao}' + ' trước ngày... tháng ... năm ...'
$new1 = 'nêu trên với ${tenLanhDao} trước ngày... tháng ... năm ...'

$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false,
                                   $true, 1, $false, $new1, 2)
Write-Host "replace1:" $found1

# --- Sentence 2 (paragraph starting with Dieu 2. Nguoi dung dau...) ---
# Originally split across 3 runs:
#   "Người đứng đầu ${tenCoQuanDuocGiaoNhiemVuXM}, ....…………………..." | "(3)" | " chịu trách nhiệm thi hành Quyết định này."
# Target: merge into a single run with the full text unchanged.
$old2 = 'Người đứng đầu ${tenCoQuanDuocGiaoNhiemVuXM}, ....…………………...' + '(3)' + ' chịu trách nhiệm thi hành Quyết định này.'
$new2 = 'Người đứng đầu ${tenCoQuanDuocGiaoNhiemVuXM}, ....…………………...(3) chịu trách nhiệm thi hành Quyết định này.'

$found2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false,
                                   $true, 1, $false, $new2, 2)
Write-Host "replace2:" $found2
